# The commit just swaps the `name` attribute values baked into the
# wp:docPr / pic:cNvPr elements of the four inline Pearson/BTEC logo
# pictures that live in the document's headers and footers:
#   PearsonLogo pictures (descr=...PearsonLogo.png): image1.png -> image2.png
#   BTEC logo pictures   (descr=BTec_Logo-Orange):   image2.jpg -> image1.jpg
#
# InlineShape objects have no settable .Name in the Word object model
# (that's the cNvPr/docPr "name" baked into the drawing XML, which Word's
# automation surface never exposes for inline pictures), so the only way
# to reach it through COM is the document's raw WordOpenXML round trip.

$d = $word.ActiveDocument

$xml = $d.WordOpenXML

# Pearson Edexcel logo: id="4" (default footer) and id="2" (first-page footer)
$xml = $xml.Replace(
    '<wp:docPr descr="Y:\Together Design\Pearson Edexcel PowerPoint amends\Assets\PearsonLogo.png" id="4" name="image1.png"/>',
    '<wp:docPr descr="Y:\Together Design\Pearson Edexcel PowerPoint amends\Assets\PearsonLogo.png" id="4" name="image2.png"/>')
$xml = $xml.Replace(
    '<wp:docPr descr="Y:\Together Design\Pearson Edexcel PowerPoint amends\Assets\PearsonLogo.png" id="2" name="image1.png"/>',
    '<wp:docPr descr="Y:\Together Design\Pearson Edexcel PowerPoint amends\Assets\PearsonLogo.png" id="2" name="image2.png"/>')
$xml = $xml.Replace(
    '<pic:cNvPr descr="Y:\Together Design\Pearson Edexcel PowerPoint amends\Assets\PearsonLogo.png" id="0" name="image1.png"/>',
    '<pic:cNvPr descr="Y:\Together Design\Pearson Edexcel PowerPoint amends\Assets\PearsonLogo.png" id="0" name="image2.png"/>')

# BTEC logo: id="1" (first-page header) and id="3" (default header)
$xml = $xml.Replace(
    '<wp:docPr descr="BTec_Logo-Orange" id="1" name="image2.jpg"/>',
    '<wp:docPr descr="BTec_Logo-Orange" id="1" name="image1.jpg"/>')
$xml = $xml.Replace(
    '<wp:docPr descr="BTec_Logo-Orange" id="3" name="image2.jpg"/>',
    '<wp:docPr descr="BTec_Logo-Orange" id="3" name="image1.jpg"/>')
$xml = $xml.Replace(
    '<pic:cNvPr descr="BTec_Logo-Orange" id="0" name="image2.jpg"/>',
    '<pic:cNvPr descr="BTec_Logo-Orange" id="0" name="image1.jpg"/>')

$d.WordOpenXML = $xml

Write-Output "done"
